$wb = $excel.ActiveWorkbook

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 159.28572
$ws.Cells.Item(33, 9).Value = 163.84616
$ws.Cells.Item(33, 11).Value = 163.84616
$ws.Cells.Item(33, 13).Value = 65.15384

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 5562.25
$ws.Cells.Item(43, 10).Value = 8124.5
$ws.Cells.Item(43, 12).Value = 8124.5
$ws.Cells.Item(43, 14).Value = -8262.5

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 3676.1177
$ws.Cells.Item(86, 10).Value = 3753.6155
$ws.Cells.Item(86, 12).Value = 3753.6155
$ws.Cells.Item(86, 14).Value = -5999.6155

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 3676.1177
$ws.Cells.Item(89, 10).Value = 3753.6155
$ws.Cells.Item(89, 12).Value = 18768.0775
$ws.Cells.Item(89, 14).Value = -30000.0775

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 935
$ws.Cells.Item(98, 9).Value = 746.7273
$ws.Cells.Item(98, 11).Value = 746.7273
$ws.Cells.Item(98, 13).Value = 751.2727

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 935
$ws.Cells.Item(122, 9).Value = 746.7273
$ws.Cells.Item(122, 11).Value = 2240.1819
$ws.Cells.Item(122, 13).Value = 209.8181

# ARM!row5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 992.5
$ws.Cells.Item(5, 10).Value = 985
$ws.Cells.Item(5, 12).Value = 985
$ws.Cells.Item(5, 14).Value = -1209

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1470.9678
$ws.Cells.Item(74, 9).Value = 1135.0176
$ws.Cells.Item(74, 11).Value = 1135.0176
$ws.Cells.Item(74, 13).Value = -261.0175999999999

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1470.9678
$ws.Cells.Item(77, 9).Value = 1135.0176
$ws.Cells.Item(77, 11).Value = 5675.088
$ws.Cells.Item(77, 13).Value = -1307.088

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1025
$ws.Cells.Item(88, 10).Value = 1750
$ws.Cells.Item(88, 12).Value = 1750
$ws.Cells.Item(88, 14).Value = -2562

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 1025
$ws.Cells.Item(91, 10).Value = 1750
$ws.Cells.Item(91, 12).Value = 1750
$ws.Cells.Item(91, 14).Value = -4558

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2705.9443
$ws.Cells.Item(122, 9).Value = 1655.1818
$ws.Cells.Item(122, 10).Value = 4357.143
$ws.Cells.Item(122, 11).Value = 4965.5454
$ws.Cells.Item(122, 12).Value = 13071.429
$ws.Cells.Item(122, 13).Value = -2515.5454
$ws.Cells.Item(122, 14).Value = -17971.429

# BSM!row4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 992.5
$ws.Cells.Item(4, 10).Value = 985
$ws.Cells.Item(4, 12).Value = 985
$ws.Cells.Item(4, 14).Value = -1215

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 864.13794
$ws.Cells.Item(94, 9).Value = 817.03705
$ws.Cells.Item(94, 10).Value = 1500
$ws.Cells.Item(94, 11).Value = 817.03705
$ws.Cells.Item(94, 12).Value = 1500
$ws.Cells.Item(94, 13).Value = -366.03705
$ws.Cells.Item(94, 14).Value = -2402

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 720.1667
$ws.Cells.Item(107, 9).Value = 709.64703
$ws.Cells.Item(107, 11).Value = 709.64703
$ws.Cells.Item(107, 13).Value = 1210.35297

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 675
$ws.Cells.Item(16, 9).Value = 667
$ws.Cells.Item(16, 10).Value = 695.8
$ws.Cells.Item(16, 11).Value = 667
$ws.Cells.Item(16, 12).Value = 695.8
$ws.Cells.Item(16, 13).Value = -380
$ws.Cells.Item(16, 14).Value = -1269.8

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5459.1665
$ws.Cells.Item(31, 9).Value = 2658.2856
$ws.Cells.Item(31, 10).Value = 7241.5454
$ws.Cells.Item(31, 11).Value = 2658.2856
$ws.Cells.Item(31, 12).Value = 7241.5454
$ws.Cells.Item(31, 13).Value = -2363.2856
$ws.Cells.Item(31, 14).Value = -7831.5454

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5459.1665
$ws.Cells.Item(34, 9).Value = 2658.2856
$ws.Cells.Item(34, 10).Value = 7241.5454
$ws.Cells.Item(34, 11).Value = 2658.2856
$ws.Cells.Item(34, 12).Value = 7241.5454
$ws.Cells.Item(34, 13).Value = -2456.2856
$ws.Cells.Item(34, 14).Value = -7645.5454

# CRP!row52
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 90000
$ws.Cells.Item(52, 10).Value = 90000
$ws.Cells.Item(52, 12).Value = 90000
$ws.Cells.Item(52, 14).Value = -90588

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 61839.43
$ws.Cells.Item(62, 9).Value = 5456
$ws.Cells.Item(62, 10).Value = 202798
$ws.Cells.Item(62, 11).Value = 5456
$ws.Cells.Item(62, 12).Value = 202798
$ws.Cells.Item(62, 13).Value = -4832
$ws.Cells.Item(62, 14).Value = -204046

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 61839.43
$ws.Cells.Item(65, 9).Value = 5456
$ws.Cells.Item(65, 10).Value = 202798
$ws.Cells.Item(65, 11).Value = 27280
$ws.Cells.Item(65, 12).Value = 1013990
$ws.Cells.Item(65, 13).Value = -24160
$ws.Cells.Item(65, 14).Value = -1020230

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 813
$ws.Cells.Item(107, 9).Value = 905
$ws.Cells.Item(107, 10).Value = 721
$ws.Cells.Item(107, 11).Value = 905
$ws.Cells.Item(107, 12).Value = 721
$ws.Cells.Item(107, 13).Value = 1015
$ws.Cells.Item(107, 14).Value = -4561

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 675
$ws.Cells.Item(113, 9).Value = 667
$ws.Cells.Item(113, 10).Value = 695.8
$ws.Cells.Item(113, 11).Value = 667
$ws.Cells.Item(113, 12).Value = 695.8
$ws.Cells.Item(113, 13).Value = 1503
$ws.Cells.Item(113, 14).Value = -5035.8

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2132.4614
$ws.Cells.Item(122, 9).Value = 2190.4348
$ws.Cells.Item(122, 11).Value = 6571.3044
$ws.Cells.Item(122, 13).Value = -4121.3044

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2885.625
$ws.Cells.Item(132, 9).Value = 2885.625
$ws.Cells.Item(132, 11).Value = 8656.875
$ws.Cells.Item(132, 13).Value = -6126.875

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1959.4054
$ws.Cells.Item(134, 9).Value = 1464.5555
$ws.Cells.Item(134, 11).Value = 4393.666499999999
$ws.Cells.Item(134, 13).Value = -1858.666499999999

# CUL!row46
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 10000000
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 13).ClearContents()

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7122.5557
$ws.Cells.Item(70, 9).Value = 5806.5
$ws.Cells.Item(70, 10).Value = 7498.5713
$ws.Cells.Item(70, 11).Value = 5806.5
$ws.Cells.Item(70, 12).Value = 7498.5713
$ws.Cells.Item(70, 13).Value = -5536.5
$ws.Cells.Item(70, 14).Value = -8038.5713

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 7122.5557
$ws.Cells.Item(73, 9).Value = 5806.5
$ws.Cells.Item(73, 10).Value = 7498.5713
$ws.Cells.Item(73, 11).Value = 5806.5
$ws.Cells.Item(73, 12).Value = 7498.5713
$ws.Cells.Item(73, 13).Value = -4870.5
$ws.Cells.Item(73, 14).Value = -9370.5713

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3553.25
$ws.Cells.Item(113, 9).Value = 3947.5
$ws.Cells.Item(113, 10).Value = 3421.8333
$ws.Cells.Item(113, 11).Value = 3947.5
$ws.Cells.Item(113, 12).Value = 3421.8333
$ws.Cells.Item(113, 13).Value = -1777.5
$ws.Cells.Item(113, 14).Value = -7761.8333

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 48775.816
$ws.Cells.Item(122, 9).Value = 2643.75
$ws.Cells.Item(122, 10).Value = 104134.3
$ws.Cells.Item(122, 11).Value = 7931.25
$ws.Cells.Item(122, 12).Value = 312402.9
$ws.Cells.Item(122, 13).Value = -5481.25
$ws.Cells.Item(122, 14).Value = -317302.9

# LTW!row3
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 348
$ws.Cells.Item(3, 9).Value = 22
$ws.Cells.Item(3, 11).Value = 22
$ws.Cells.Item(3, 13).Value = 90

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2488.8333
$ws.Cells.Item(7, 9).Value = 2486.6
$ws.Cells.Item(7, 11).Value = 2486.6
$ws.Cells.Item(7, 13).Value = -2374.6

# LTW!row15
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(15, 8).Value = 348
$ws.Cells.Item(15, 9).Value = 22
$ws.Cells.Item(15, 11).Value = 22
$ws.Cells.Item(15, 13).Value = 148

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6822
$ws.Cells.Item(22, 9).Value = 1524.75
$ws.Cells.Item(22, 10).Value = 11059.8
$ws.Cells.Item(22, 11).Value = 1524.75
$ws.Cells.Item(22, 12).Value = 11059.8
$ws.Cells.Item(22, 13).Value = -1229.75
$ws.Cells.Item(22, 14).Value = -11649.8

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 6822
$ws.Cells.Item(27, 9).Value = 1524.75
$ws.Cells.Item(27, 10).Value = 11059.8
$ws.Cells.Item(27, 11).Value = 1524.75
$ws.Cells.Item(27, 12).Value = 11059.8
$ws.Cells.Item(27, 13).Value = -1417.75
$ws.Cells.Item(27, 14).Value = -11273.8

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 244.4
$ws.Cells.Item(55, 9).Value = 255.1579
$ws.Cells.Item(55, 10).Value = 40
$ws.Cells.Item(55, 11).Value = 255.1579
$ws.Cells.Item(55, 12).Value = 40
$ws.Cells.Item(55, 13).Value = -82.15790000000001
$ws.Cells.Item(55, 14).Value = -386

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 10502
$ws.Cells.Item(61, 9).Value = 15999
$ws.Cells.Item(61, 11).Value = 15999
$ws.Cells.Item(61, 13).Value = -15797

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1999
$ws.Cells.Item(68, 9).Value = 1999
$ws.Cells.Item(68, 11).Value = 1999
$ws.Cells.Item(68, 13).Value = -1250

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 1999
$ws.Cells.Item(71, 9).Value = 1999
$ws.Cells.Item(71, 11).Value = 9995
$ws.Cells.Item(71, 13).Value = -6251

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2410.2122
$ws.Cells.Item(82, 9).Value = 2664.3333
$ws.Cells.Item(82, 10).Value = 1732.5555
$ws.Cells.Item(82, 11).Value = 2664.3333
$ws.Cells.Item(82, 12).Value = 1732.5555
$ws.Cells.Item(82, 13).Value = -2303.3333
$ws.Cells.Item(82, 14).Value = -2454.5555

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2410.2122
$ws.Cells.Item(85, 9).Value = 2664.3333
$ws.Cells.Item(85, 10).Value = 1732.5555
$ws.Cells.Item(85, 11).Value = 2664.3333
$ws.Cells.Item(85, 12).Value = 1732.5555
$ws.Cells.Item(85, 13).Value = -1416.3333
$ws.Cells.Item(85, 14).Value = -4228.5555

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 10502
$ws.Cells.Item(113, 9).Value = 15999
$ws.Cells.Item(113, 11).Value = 15999
$ws.Cells.Item(113, 13).Value = -13829

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 2488.8333
$ws.Cells.Item(126, 9).Value = 2486.6
$ws.Cells.Item(126, 11).Value = 7459.799999999999
$ws.Cells.Item(126, 13).Value = -4989.799999999999

# WVR!row26
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 352083.34
$ws.Cells.Item(26, 10).Value = 352083.34
$ws.Cells.Item(26, 12).Value = 352083.34
$ws.Cells.Item(26, 14).Value = -352669.34

# WVR!row43
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 22500
$ws.Cells.Item(43, 10).Value = 22500
$ws.Cells.Item(43, 12).Value = 22500
$ws.Cells.Item(43, 14).Value = -22798

# WVR!row80
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 90000
$ws.Cells.Item(80, 10).Value = 90000
$ws.Cells.Item(80, 12).Value = 90000
$ws.Cells.Item(80, 14).Value = -91996

# WVR!row83
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(83, 8).Value = 90000
$ws.Cells.Item(83, 10).Value = 90000
$ws.Cells.Item(83, 12).Value = 270000
$ws.Cells.Item(83, 14).Value = -279984

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1532.0952
$ws.Cells.Item(136, 9).Value = 1167.2632
$ws.Cells.Item(136, 11).Value = 3501.7896
$ws.Cells.Item(136, 13).Value = -951.7896000000001
